$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must remain text
# (matching the source data format), so we force a text number format,
# assign the value, then restore the default style (no explicit number format).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.129'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.402'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.135'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '436.29'
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.78'
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.75'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '162.83'
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.42'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.794'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '321.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.976'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.16'
$ws.Range("D51").Style = "Normal"

# Remaining cells can be assigned directly
$ws.Range("D2").Value = '66.536.36'
$ws.Range("E2").Value = '  -4.61%  '
$ws.Range("D3").Value = '3.332.13'
$ws.Range("E3").Value = '  -1.57%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  -3.41%  '
$ws.Range("E6").Value = '  -5.50%  '
$ws.Range("E7").Value = '  +2.69%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -3.61%  '
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("E11").Value = '  -4.03%  '
$ws.Range("D12").Value = '3.909.97'
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("E14").Value = '  -6.25%  '
$ws.Range("D15").Value = '66.642.65'
$ws.Range("E15").Value = '  -4.40%  '
$ws.Range("E16").Value = '  -3.03%  '
$ws.Range("D17").Value = '3.317.44'
$ws.Range("E17").Value = '  -2.42%  '
$ws.Range("E18").Value = '  -3.53%  '
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("E21").Value = '  -3.03%  '
$ws.Range("E22").Value = '  -3.48%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  -4.35%  '
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("E27").Value = '  -5.63%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("E29").Value = '  -3.37%  '
$ws.Range("E30").Value = '  -2.88%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  -6.48%  '
$ws.Range("E33").Value = '  -4.42%  '
$ws.Range("E34").Value = '  -3.61%  '
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("E36").Value = '  -5.98%  '
$ws.Range("E37").Value = '  -1.79%  '
$ws.Range("E38").Value = '  -8.21%  '
$ws.Range("D39").Value = '2.828.91'
$ws.Range("E40").Value = '  -2.48%  '
$ws.Range("E41").Value = '  -4.29%  '
$ws.Range("E42").Value = '  -5.74%  '
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("E44").Value = '  -3.29%  '
$ws.Range("E45").Value = '  -4.55%  '
$ws.Range("E46").Value = '  -7.01%  '
$ws.Range("E47").Value = '  -5.57%  '
$ws.Range("E48").Value = '  -4.07%  '
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("E50").Value = '  -4.63%  '
$ws.Range("E51").Value = '  -2.82%  '
